$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 650 (existing rows 650-685 shift down to 652-687)
$ws.Rows.Item(650).Resize(2).Insert()

# New row 650: Primera, fecha 44706
$ws.Cells.Item(650,1).Value  = 3
$ws.Cells.Item(650,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(650,3).Value  = "Coquimbo"
$ws.Cells.Item(650,4).Value  = 44706
$ws.Cells.Item(650,5).Value  = 5
$ws.Cells.Item(650,6).Value  = 100112023
$ws.Cells.Item(650,7).Value  = "Brócoli"
$ws.Cells.Item(650,8).Value  = "Sin especificar"
$ws.Cells.Item(650,9).Value  = "Primera"
$ws.Cells.Item(650,10).Value = 950
$ws.Cells.Item(650,11).Value = 900
$ws.Cells.Item(650,12).Value = 900
$ws.Cells.Item(650,13).Value = 900
$ws.Cells.Item(650,14).Value = "$/unidad"
$ws.Cells.Item(650,15).Value = "Provincia de Quillota"
$ws.Cells.Item(650,16).Value = 900
$ws.Cells.Item(650,17).Value = 1
$ws.Cells.Item(650,18).Value = "Hortaliza"

# New row 651: Segunda, fecha 44706
$ws.Cells.Item(651,1).Value  = 3
$ws.Cells.Item(651,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(651,3).Value  = "Coquimbo"
$ws.Cells.Item(651,4).Value  = 44706
$ws.Cells.Item(651,5).Value  = 5
$ws.Cells.Item(651,6).Value  = 100112023
$ws.Cells.Item(651,7).Value  = "Brócoli"
$ws.Cells.Item(651,8).Value  = "Sin especificar"
$ws.Cells.Item(651,9).Value  = "Segunda"
$ws.Cells.Item(651,10).Value = 1780
$ws.Cells.Item(651,11).Value = 600
$ws.Cells.Item(651,12).Value = 700
$ws.Cells.Item(651,13).Value = 651
$ws.Cells.Item(651,14).Value = "$/unidad"
$ws.Cells.Item(651,15).Value = "Provincia de Quillota"
$ws.Cells.Item(651,16).Value = 651
$ws.Cells.Item(651,17).Value = 1
$ws.Cells.Item(651,18).Value = "Hortaliza"
